# Meeting diary update: add the Q4-completion / Q5 kickoff meeting entry
# as row 13, mirroring the format of the preceding meeting rows (7-12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New meeting record values (Date, Time start, Time end, Members present, Discussions)
$ws.Range("A13").Value = 45203
$ws.Range("B13").Value = 0.91666666666666663
$ws.Range("C13").Value = 0.99305555555555547
$ws.Range("D13").Value = "All"
$ws.Range("E13").Value = "- Done with Q4`n- Merging of Q4 to be done on 04/10/2023`n- Q5 coding part done by Evan`n- Deadline for Q5 explanation is 05/10/2023`n- Deadline for Q5 VIA is 06/10/2023"

# Copy formatting from the row above so styles (date/time/centered/wrap) match
$ws.Range("A12:E12").Copy()
$ws.Range("A13:E13").PasteSpecial(-4122)

# Match the row height used for this wrapped discussion entry
$ws.Rows(13).RowHeight = 78

# Move selection to reflect the newly added row being the focus, as in the
# saved workbook (cursor parked just below the new entry)
[void]$ws.Range("E14").Select()
